# Disable Get OCR Skills Gap - Not working
# Reset the "Score (%)" column values (which previously showed OCR-derived
# skills-match percentages) back to the default "0% skills matched" value,
# since the OCR skills gap feature is being disabled / not working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = "0% skills matched"
$ws.Range("I4").Value = "0% skills matched"
$ws.Range("I5").Value = "0% skills matched"
